# Atlas UAT Test Full Integration Template - "Add files via upload" edit
#
# Renames the "Baggage" concept to "Ancillary Baggage" throughout the sample
# data sheet, re-sequences the Sequence column (1..5), widens column I,
# doubles the height of row 4 (now-longer wrapped text), moves the selection
# / scroll position, and adds a reviewer comment on J1 ("Samples").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Sequence"
$ws.Range("B1").Value = "PassengerType"
$ws.Range("C1").Value = "TripType"
$ws.Range("D1").Value = "Direct/Connection"
$ws.Range("E1").Value = "PaymentMothod"
$ws.Range("G1").Value = "Post Ticketing Purchase"
$ws.Range("H1").Value = "Refund"
$ws.Range("I1").Value = "Schedule Change Notification"
$ws.Range("J1").Value = "Samples"
$ws.Range("K1").Value = "Airline PNR"

# ---------------------------------------------------------------------
# Sequence numbers 1..5 (rows 2..6), and the passenger/trip/payment
# columns (B-E), which keep their original wording.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "1ADT"
$ws.Range("C2").Value = "Oneway"
$ws.Range("D2").Value = "Connection Flights"
$ws.Range("E2").Value = "Prepayment"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "1ADT"
$ws.Range("C3").Value = "Oneway"
$ws.Range("D3").Value = "Connection Flights"
$ws.Range("E3").Value = "Prepayment"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "2ADT1CHD"
$ws.Range("C4").Value = "Roundtrip"
$ws.Range("D4").Value = "Direct Flights"
$ws.Range("E4").Value = "Prepayment"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "2ADT1CHD"
$ws.Range("C5").Value = "Oneway"
$ws.Range("D5").Value = "Direct Flights"
$ws.Range("E5").Value = "Clients' VCC"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "none"
$ws.Range("C6").Value = "none"
$ws.Range("D6").Value = "none"
$ws.Range("E6").Value = "none"

# ---------------------------------------------------------------------
# Column F (Baggage -> Ancillary Baggage), edited top to bottom.
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Ancillary Baggage"
$ws.Range("F2").Value = "Without Ancillary Baggage"
$ws.Range("F3").Value = "With Ancillary Baggage"
$ws.Range("F4").Value = "1ADT with Ancillary Baggage in both bounds"
$ws.Range("F5").Value = "All of them with Ancillary Baggage"
$ws.Range("F6").Value = "none"

# ---------------------------------------------------------------------
# Column G (Extra Baggage -> Extra Ancillary Baggage), top to bottom.
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "Extra Ancillary Baggage"
$ws.Range("G3").Value = "Extra Ancillary Baggage"
$ws.Range("G4").Value = "1ADT extra Ancillary Baggage in outbound" + $NL + "1CHD extra Ancillary Baggage in inbound"
$ws.Range("G5").Value = "1ADT extra Ancillary baggage"
$ws.Range("G6").Value = "none"
$ws.Rows(4).RowHeight = 58

# ---------------------------------------------------------------------
# Column H (Refund detail text, unchanged wording).
# ---------------------------------------------------------------------
$ws.Range("H2").Value = "Total refund"
$ws.Range("H3").Value = "Total refund"
$ws.Range("H4").Value = "1ADT1CHD inbound refund"
$ws.Range("H5").Value = "1ADT1CHD refund"
$ws.Range("H6").Value = "none"

# ---------------------------------------------------------------------
# Column I (Schedule Change Notification detail, unchanged wording).
# ---------------------------------------------------------------------
$ws.Range("I2").Value = "none"
$ws.Range("I3").Value = "none"
$ws.Range("I4").Value = "none"
$ws.Range("I5").Value = "none"
$ws.Range("I6").Value = "Register your schedule change notification receiver"

# ---------------------------------------------------------------------
# Column J (Samples city pairs, unchanged wording).
# ---------------------------------------------------------------------
$ws.Range("J2").Value = "JKT-SUB  2-20"
$ws.Range("J3").Value = "JKT-SUB  2-20"
$ws.Range("J4").Value = "MNL-CEB  2-20/2-24"
$ws.Range("J5").Value = "DVO-CEB  2-20"
$ws.Range("J6").Value = "none"

# ---------------------------------------------------------------------
# L1 header: orderNo -> OrderNo (edited last).
# ---------------------------------------------------------------------
$ws.Range("L1").Value = "OrderNo"

# ---------------------------------------------------------------------
# Column I (9) widened to fit the longer "Schedule Change Notification"
# style header / content.
# ---------------------------------------------------------------------
$ws.Columns(9).ColumnWidth = 42.83

# ---------------------------------------------------------------------
# View state: scrolled one column right, selection parked at J7.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J7").Select() | Out-Null

# ---------------------------------------------------------------------
# Reviewer note on the "Samples" header (J1).
# ---------------------------------------------------------------------
$commentText = "Behram Kotwal:" + $NL + "These are just samples. Any city pair and dates can be used which fulfill the user case."
$ws.Range("J1").AddComment($commentText) | Out-Null
